$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.485.78'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.749.26'
$ws.Range("E3").Value = '  -1.34%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").Value = '324.08'
$ws.Range("E5").Value = '  +0.69%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '0.4456'
$ws.Range("E7").Value = '  +4.48%  '
$ws.Range("D8").Value = '0.3572'
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").Value = '0.07482'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").Value = '42.12'
$ws.Range("E10").Value = '  -4.72%  '
$ws.Range("D11").Value = '1.088'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("E12").Value = '  +0.71%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '20.70'
$c.ClearFormats()
$ws.Range("E13").Value = '  -3.99%  '
$ws.Range("D14").Value = '6.011'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.760.64'
$ws.Range("E15").Value = '  -1.92%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '7.091'
$ws.Range("E16").Value = '  -2.83%  '
$ws.Range("D17").Value = '92.85'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").Value = '0.00001059'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '0.06405'
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '16.79'
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").Value = '5.803'
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").Value = '27.545.53'
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("D24").Value = '11.16'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").Value = '2.108'
$ws.Range("E25").Value = '  -2.63%  '
$ws.Range("D26").Value = '162.56'
$ws.Range("D27").Value = '20.44'
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("D28").Value = '1.948.73'
$ws.Range("E28").Value = '  -2.20%  '
$ws.Range("D29").Value = '2.071'
$ws.Range("E29").Value = '  -4.30%  '
$ws.Range("D30").Value = '125.71'
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '1.072'
$ws.Range("E31").Value = '  -8.05%  '
$ws.Range("D32").Value = '3.667'
$ws.Range("E32").Value = '  +4.54%  '
$ws.Range("D33").Value = '0.09046'
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '5.498'
$ws.Range("E34").Value = '  -3.37%  '
$ws.Range("D35").Value = '11.92'
$ws.Range("E35").Value = '  -5.75%  '
$ws.Range("D36").Value = '0.02282'
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("D37").Value = '0.2094'
$ws.Range("E37").Value = '  -0.80%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.06010'
$c.ClearFormats()
$ws.Range("E38").Value = '  -0.70%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.6350'
$c.ClearFormats()
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").Value = '4.928'
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("D41").Value = '1.203'
$ws.Range("E41").Value = '  +2.12%  '
$ws.Range("D42").Value = '1.376'
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").Value = '7.736'
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("D44").Value = '13.15'
$ws.Range("E44").Value = '  -3.57%  '
$ws.Range("D45").Value = '3.714'
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("D46").Value = '0.5885'
$ws.Range("E46").Value = '  -1.31%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '121.60'
$c.ClearFormats()
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").Value = '1.143'
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").Value = '0.06831'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").Value = '72.25'
$ws.Range("E51").Value = '  -2.99%  '
